$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the stat data (columns C through DK) between rows 7 and 8,
# leaving League (A) and Team (B) columns untouched since both
# players belong to Serie_A / Bologna. Column DL (goalsPrevented) is
# blank for both rows, so it is left alone.
$row7 = $ws.Range("C7:DK7").Value2
$row8 = $ws.Range("C8:DK8").Value2

$ws.Range("C7:DK7").Value2 = $row8
$ws.Range("C8:DK8").Value2 = $row7
